# Update the acquisition datetime (column A, "取得日時") for all data rows
# on the "ランサーズ" sheet from "2025-09-13 12:40:45" to "2025-09-13 18:21:56".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-13 12:40:45"
$newValue = "2025-09-13 18:21:56"

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
